# Update faturamento_diario_lojas data:
# Column F (day 5) values were missing (0) and are now filled in for each store,
# which also changes the row total in column AG.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bibi Cell Mundi
$ws.Range("F2").Value = 14123.07
$ws.Range("AG2").Value = 39866.69

# Row 3 - Bibi Cell Vieiralves
$ws.Range("F3").Value = 3535.9
$ws.Range("AG3").Value = 22931.3

# Row 4 - Bibi Cell Manauara
$ws.Range("F4").Value = 1843
$ws.Range("AG4").Value = 16721.75

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("F5").Value = 1546
$ws.Range("AG5").Value = 13944.55

# Row 6 - total
$ws.Range("F6").Value = 21047.97
$ws.Range("AG6").Value = 93464.28999999999
